$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "1.001"); force text
# formatting before assignment so Excel does not coerce it to a Double, then
# restore the default "Normal" style so no stray number-format style lingers
# on the cell (matches the source workbook, where these cells carry no style).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "24.780.54"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.699.89"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  -0.31%  "

# Row 5
Set-TextValue $ws.Range("D5") "314.68"

# Row 6
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.3985"
$ws.Range("E7").Value = "  +2.57%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.4034"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D9") "1.002"
$ws.Range("E9").Value = "  -0.14%  "

# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D10") "1.471"
$ws.Range("E10").Value = "  -1.55%  "

# Row 11
Set-TextValue $ws.Range("D11") "53.44"
$ws.Range("E11").Value = "  +1.73%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.08799"
$ws.Range("E12").Value = "  +0.49%  "

# Row 13
Set-TextValue $ws.Range("D13") "26.11"
$ws.Range("E13").Value = "  +3.78%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.546"
$ws.Range("E14").Value = "  +0.27%  "

# Row 15
Set-TextValue $ws.Range("D15") "7.986"
$ws.Range("E15").Value = "  -0.02%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.00001347"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.697.68"
$ws.Range("E17").Value = "  +0.65%  "

# Row 18
Set-TextValue $ws.Range("D18") "95.72"
$ws.Range("E18").Value = "  -2.75%  "

# Row 19
Set-TextValue $ws.Range("D19") "0.07176"
$ws.Range("E19").Value = "  +1.10%  "

# Row 20
Set-TextValue $ws.Range("D20") "20.89"
$ws.Range("E20").Value = "  +4.17%  "

# Row 21
Set-TextValue $ws.Range("D21") "7.349"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22
Set-TextValue $ws.Range("D22") "1.000"
$ws.Range("E22").Value = "  -0.32%  "

# Row 23
Set-TextValue $ws.Range("D23") "14.41"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
Set-TextValue $ws.Range("D24") "24.778.02"
$ws.Range("E24").Value = "  +1.20%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.367"
$ws.Range("E25").Value = "  +0.75%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.926"
$ws.Range("E26").Value = "  -1.52%  "

# Row 27
$ws.Range("E27").Value = "  +1.93%  "

# Row 28
Set-TextValue $ws.Range("D28") "6.165"
$ws.Range("E28").Value = "  +18.04%  "

# Row 29
Set-TextValue $ws.Range("D29") "161.77"
$ws.Range("E29").Value = "  -0.42%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D30") "8.478"
$ws.Range("E30").Value = "  -3.35%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D31") "144.10"
$ws.Range("E31").Value = "  +5.15%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.380"
$ws.Range("E32").Value = "  +21.33%  "

# Row 33
Set-TextValue $ws.Range("D33") "1.897.85"
$ws.Range("E33").Value = "  +1.36%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.08637"
$ws.Range("E34").Value = "  -2.37%  "

# Row 35
Set-TextValue $ws.Range("D35") "7.341"
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.03172"
$ws.Range("E36").Value = "  +8.47%  "

# Row 37
Set-TextValue $ws.Range("D37") "1.037"
$ws.Range("E37").Value = "  +0.05%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.2837"
$ws.Range("E38").Value = "  +0.43%  "

# Row 39
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D39") "10.78"
$ws.Range("E39").Value = "  -0.25%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D40") "0.09451"
$ws.Range("E40").Value = "  +3.53%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.8298"
$ws.Range("E41").Value = "  +4.30%  "

# Row 42
Set-TextValue $ws.Range("D42") "14.22"
$ws.Range("E42").Value = "  -0.04%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.476"
$ws.Range("E43").Value = "  +1.50%  "

# Row 44
Set-TextValue $ws.Range("D44") "17.77"
$ws.Range("E44").Value = "  +7.16%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.705"
$ws.Range("E45").Value = "  +3.26%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.7426"
$ws.Range("E46").Value = "  +2.21%  "

# Row 47
Set-TextValue $ws.Range("D47") "4.211"
$ws.Range("E47").Value = "  +0.10%  "

# Row 48
Set-TextValue $ws.Range("D48") "1.388"
$ws.Range("E48").Value = "  +2.60%  "

# Row 49
Set-TextValue $ws.Range("D49") "1.000"
$ws.Range("E49").Value = "  -0.25%  "

# Row 50
Set-TextValue $ws.Range("D50") "0.08376"
$ws.Range("E50").Value = "  +4.44%  "

# Row 51
Set-TextValue $ws.Range("D51") "139.73"
$ws.Range("E51").Value = "  +0.87%  "
